$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimParameters")

# Explored different data generating parameters: drop the "Low" severity
# relative-risk multipliers for both Abortion and Preeclampsia from 0.25 to 0.1.
$ws.Range("B4").Value = 0.1
$ws.Range("B8").Value = 0.1

# Reflect the author's last selection on the SimParameters sheet.
$ws.Activate()
$ws.Range("B9").Select()
